$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.215.38'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.888.19'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.36%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '482.60'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.62'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.87%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.64%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.998'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.744'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.11%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000356'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.88%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.10'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.86%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.52'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.505.17'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.41%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.962.99'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.26'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -2.09%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.96'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.47%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.70%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '68.248.46'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.31%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '429.51'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.59'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +8.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.81'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.30%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.58'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +17.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '88.73'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +2.28%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.02'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.08%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.62%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '719.63'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.50'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.24%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.62%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.14%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '61.75'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +6.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.07'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +10.77%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0₃0875'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '40.90'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.402'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +17.59%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.12%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.11%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +6.37%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +6.78%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.95'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.34%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.76%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.37'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +5.88%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0362'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +32.49%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.72%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.82%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '144.53'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.49%  '
